$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab) to reflect the new "through" date
$ws.Name = "Through 2021-12-20"

# Row 14 - "December (through 12-20)" data
$ws.Range("A14").Value = "December (through 12-20)"

$ws.Range("C14").Value = 24
$ws.Range("D14").Value = 0.1111

$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 58
$ws.Range("G14").Value = 0.0938

$ws.Range("H14").Value = 9
$ws.Range("I14").Value = 68
$ws.Range("J14").Value = 0.1169

$ws.Range("L14").Value = 40
$ws.Range("M14").Value = 0.0909

$ws.Range("O14").Value = 30
$ws.Range("P14").Value = 0.0909

$ws.Range("R14").Value = 89
$ws.Range("S14").Value = 0.043

$ws.Range("U14").Value = 136
$ws.Range("V14").Value = 0.0145

# Row 15 - "Total" data
$ws.Range("C15").Value = 282
$ws.Range("D15").Value = 0.1132

$ws.Range("E15").Value = 65
$ws.Range("F15").Value = 562
$ws.Range("G15").Value = 0.1037

$ws.Range("H15").Value = 72
$ws.Range("I15").Value = 826
$ws.Range("J15").Value = 0.0802

$ws.Range("L15").Value = 648
$ws.Range("M15").Value = 0.1074

$ws.Range("O15").Value = 510
$ws.Range("P15").Value = 0.1005

$ws.Range("R15").Value = 1289
$ws.Range("S15").Value = 0.0501

$ws.Range("U15").Value = 1678
$ws.Range("V15").Value = 0.0578
